$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each test case block in the sheet starts with a header row that has
# Test Case / Expected Behavior columns (A/B) and now gets a new "Passed"
# status in column C, shown in green text.
$headerRows = @(1, 9, 15, 21, 27, 33, 39, 45)

foreach ($r in $headerRows) {
    $cell = $ws.Range("C$r")
    $cell.Value = "Passed"
    $cell.Font.Color = 5287936   # RGB(0,176,80) -> 0x00B050 green
}

# Page setup: use Letter... actually use paper size 9 (A4) and portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update the active selection/view to cell C45 (no need to scroll since it
# is already within the visible area).
$ws.Range("C45").Select()
